# BD.xlsx - actualizar fila de datos existente (Metodo CrearExcel: No agrega mas de uno)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRow = $ws.Range("A2:I2")

# Forzar formato de texto para que los valores numericos (ej. "4.0", "7.0")
# se guarden como texto, igual que en el archivo original (t="s"), y no se
# conviertan automaticamente a numeros.
$dataRow.NumberFormat = "@"

$ws.Range("A2").Value = "24/02/2020"
$ws.Range("B2").Value = "02:02:08"
$ws.Range("C2").Value = "4.0"
$ws.Range("D2").Value = "SUBDUCCION_PLACA"
$ws.Range("E2").Value = "hola"
$ws.Range("F2").Value = "7.0"
$ws.Range("G2").Value = "9.0"
$ws.Range("H2").Value = "4.0"
$ws.Range("I2").Value = "SAN_JOSE, adios"

# Restaurar el formato por defecto de las celdas (sin estilo de texto forzado)
# una vez que el tipo de celda como texto ya quedo fijado.
$dataRow.ClearFormats()
